$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Col8a1", "Itga2", "ECs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 2, 0.6666666666666666, 1.242496666666667, 3.72749, 0.3272238221337332, 0.3272238221337332, 5.682122802835555, 51.13910522552, 0.02028566996228753, 0.02028566996228752),
  @("ECs", "Col8a1", "Itga2", "FAPs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 3, 1, 1.744414, 5.233242000000001, 0.459408730644692, 0.459408730644692, 7.977465721157333, 71.79719149041601, 0.02848024274908357, 0.02848024274908357),
  @("ECs", "Col8a1", "Itga2", "sCs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 3, 1, 0.8101743333333333, 2.430523, 0.2133674472215748, 0.2133674472215748, 3.705048212367111, 33.345433911304, 0.01322734263908125, 0.01322734263908125),
  @("FAPs", "Col8a1", "Itga2", "ECs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 2, 0.6666666666666666, 1.242496666666667, 3.72749, 0.3272238221337332, 0.3272238221337332, 71.26789084127556, 641.41101757148, 0.2544325356349201, 0.2544325356349201),
  @("FAPs", "Col8a1", "Itga2", "FAPs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 3, 1, 1.744414, 5.233242000000001, 0.459408730644692, 0.459408730644692, 100.0571750969093, 900.5145758721841, 0.3572127709668332, 0.3572127709668332),
  @("FAPs", "Col8a1", "Itga2", "sCs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 3, 1, 0.8101743333333333, 2.430523, 0.2133674472215748, 0.2133674472215748, 46.47047955895511, 418.234316030596, 0.1659036321516605, 0.1659036321516605),
  @("sCs", "Col8a1", "Itga2", "ECs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 2, 0.6666666666666666, 1.242496666666667, 3.72749, 0.3272238221337332, 0.3272238221337332, 14.70709922589555, 132.36389303306, 0.05250561653652558, 0.05250561653652557),
  @("sCs", "Col8a1", "Itga2", "FAPs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 3, 1, 1.744414, 5.233242000000001, 0.459408730644692, 0.459408730644692, 20.64815985210533, 185.833438668948, 0.0737157169287752, 0.07371571692877518),
  @("sCs", "Col8a1", "Itga2", "sCs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 3, 1, 0.8101743333333333, 2.430523, 0.2133674472215748, 0.2133674472215748, 9.58981591682911, 86.30834325146199, 0.03423647243083302, 0.03423647243083301),
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
